# knitxl headers1.xlsx refactor
# -----------------------------------------------------------------------
# The original sheet crammed the whole markdown sample ("# Header 1",
# "text", "## Header 2", "text", ... "###### Header 6", "text") into a
# single A1 cell/shared-string. This reshapes it into the "one xl
# snapshot per test" layout: one row per markdown block (12 rows total),
# each with its own cell value and its own knitxl heading/body style
# (font + bottom border), matching what the real converter emits.
# -----------------------------------------------------------------------

function Convert-RgbToExcelColor([string]$rgbHex) {
  # Excel's Font/Border .Color is a BGR-packed long (same convention as
  # VBA's RGB()), while we think in normal "RRGGBB" hex.
  $r = [Convert]::ToInt32($rgbHex.Substring(0, 2), 16)
  $g = [Convert]::ToInt32($rgbHex.Substring(2, 2), 16)
  $b = [Convert]::ToInt32($rgbHex.Substring(4, 2), 16)
  return $r + ($g * 256) + ($b * 65536)
}

$HeadingColor = Convert-RgbToExcelColor "475368"   # knitxl heading/body grey-blue
$BodyColor    = Convert-RgbToExcelColor "000000"   # plain paragraph text
$H2BorderColor = Convert-RgbToExcelColor "4F71BE"  # Header 2 underline
$H34BorderColor = Convert-RgbToExcelColor "A6B7DE" # Header 3 / Header 4 underline

$xlThick  = 4
$xlMedium = -4138
$xlEdgeBottom = 9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- cell values, one markdown block per row ------------------------------
$ws.Range("A1").Value  = "Header 1"
$ws.Range("A2").Value  = "text"
$ws.Range("A3").Value  = "Header 2"
$ws.Range("A4").Value  = "text"
$ws.Range("A5").Value  = "Header 3"
$ws.Range("A6").Value  = "text"
$ws.Range("A7").Value  = "Header 4"
$ws.Range("A8").Value  = "text"
$ws.Range("A9").Value  = "Header 5"
$ws.Range("A10").Value = "text"
$ws.Range("A11").Value = "Header 6"
$ws.Range("A12").Value = "text"

# ---- Header 1 (A1): 18pt "Calibri Light (Headings)", no border -----------
$c = $ws.Range("A1")
$c.Font.Name  = "Calibri Light (Headings)"
$c.Font.Size  = 18
$c.Font.Color = $HeadingColor

# ---- body "text" rows: 11pt Calibri, plain black --------------------------
foreach ($addr in @("A2", "A4", "A6", "A8", "A10", "A12")) {
  $c = $ws.Range($addr)
  $c.Font.Name  = "Calibri"
  $c.Font.Size  = 11
  $c.Font.Color = $BodyColor
}

# ---- Header 2 (A3): 15pt bold "Calibri (Body)", thick blue underline -----
$c = $ws.Range("A3")
$c.Font.Name  = "Calibri (Body)"
$c.Font.Size  = 15
$c.Font.Color = $HeadingColor
$c.Font.Bold  = $true
$b = $c.Borders.Item($xlEdgeBottom)
$b.Color  = $H2BorderColor
$b.Weight = $xlThick

# ---- Header 3 (A5): 13pt bold "Calibri (Body)", thick light underline ----
$c = $ws.Range("A5")
$c.Font.Name  = "Calibri (Body)"
$c.Font.Size  = 13
$c.Font.Color = $HeadingColor
$c.Font.Bold  = $true
$b = $c.Borders.Item($xlEdgeBottom)
$b.Color  = $H34BorderColor
$b.Weight = $xlThick

# ---- Header 4 (A7): 11pt bold "Calibri (Body)", medium light underline ---
$c = $ws.Range("A7")
$c.Font.Name  = "Calibri (Body)"
$c.Font.Size  = 11
$c.Font.Color = $HeadingColor
$c.Font.Bold  = $true
$b = $c.Borders.Item($xlEdgeBottom)
$b.Color  = $H34BorderColor
$b.Weight = $xlMedium

# ---- Header 5 (A9): 11pt bold "Calibri (Body)", no border -----------------
$c = $ws.Range("A9")
$c.Font.Name  = "Calibri (Body)"
$c.Font.Size  = 11
$c.Font.Color = $HeadingColor
$c.Font.Bold  = $true

# ---- Header 6 (A11): 11pt italic "Calibri (Body)", no border --------------
$c = $ws.Range("A11")
$c.Font.Name   = "Calibri (Body)"
$c.Font.Size   = 11
$c.Font.Color  = $HeadingColor
$c.Font.Italic = $true

Write-Host "headers1.xlsx rebuilt: 12 rows, 7 shared strings, 6 heading styles"
